$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "57.904.84"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.958.62"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.30"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.82"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +10.97%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.515"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +4.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.954.38"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.27%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +4.30%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +2.01%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.50%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.27%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.00"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.85%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.448.64"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.35%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.99"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +8.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.956.67"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "57.865.45"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.47%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "419.56"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.61%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.52"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.707"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.57%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.41"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.72%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +4.83%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +4.07%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.02%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.71%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +7.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.54"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.65"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +4.44%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0970"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +1.93%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.73"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +7.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.957"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +6.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.07"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.48%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0709"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +14.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "8.87"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.33%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "48.31"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.71"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +15.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "387.75"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +6.87%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0348"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.39%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.107"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.721.10"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.47%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "124.74"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +5.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.239"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.00"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +3.18%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.31%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.00"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.83%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +3.54%  "

Write-Host "Applied crypto price updates"